$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.156.81"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "2.521.73"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.69"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.43"
$ws.Range("E6").Value = "  -1.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("D9").Value = "2.527.99"
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0994"
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.40"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("D14").Value = "2.965.94"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.00"
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").Value = "59.114.82"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("D18").Value = "2.525.75"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.91"
$ws.Range("E19").Value = "  -2.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.22"
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.44"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.81"
$ws.Range("E23").Value = "  +1.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.61"
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.423"
$ws.Range("E25").Value = "  -2.26%  "
$ws.Range("E26").Value = "  +1.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.77"
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("D31").Value = "0.0₃0764"
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.11"
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.13"
$ws.Range("E34").Value = "  -3.37%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.45"
$ws.Range("E35").Value = "  +2.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.52"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.19"
$ws.Range("E37").Value = "  -2.53%  "
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.98"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.63"
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.805"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.24"
$ws.Range("E42").Value = "  -6.16%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "283.16"
$ws.Range("E43").Value = "  -4.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.596"
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0927"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.22"
$ws.Range("E48").Value = "  -1.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.49"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0509"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("E51").Value = "  -1.64%  "
